$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the sample-data rows (2-7): the sheet used to repeat a single
# person's (Nguyen Duc Minh / Keppel) data with only the phone numbers and
# the picture varying per row. The "generate-list" sample now alternates
# between two people (Nguyen Van A / Nguyen Van B) sharing one email,
# company (Gooogle) and phone number, again varying only per-row phone +
# picture. ---

$rows = @(
    @{ Row=2; A=1; B="Nguyen Van A"; C="Head of IT"; D="abcd@gmail.com"; E="84 (28)12345678"; F="+84123456887"; G="https://www.google.com.vn/?hl=vi"; H="Gooogle"; I="https://picsum.photos/200/300" },
    @{ Row=3; A=2; B="Nguyen Van B"; C="Head of IT"; D="abcd@gmail.com"; E="85 (28)12345678"; F="+84123456887"; G="https://www.google.com.vn/?hl=vi"; H="Gooogle"; I="https://picsum.photos/200/301" },
    @{ Row=4; A=3; B="Nguyen Van A"; C="Head of IT"; D="abcd@gmail.com"; E="86 (28)12345678"; F="+84123456887"; G="https://www.google.com.vn/?hl=vi"; H="Gooogle"; I="https://picsum.photos/200/302" },
    @{ Row=5; A=4; B="Nguyen Van B"; C="Head of IT"; D="abcd@gmail.com"; E="87 (28)12345678"; F="+84123456887"; G="https://www.google.com.vn/?hl=vi"; H="Gooogle"; I="https://picsum.photos/200/303" },
    @{ Row=6; A=5; B="Nguyen Van A"; C="Head of IT"; D="abcd@gmail.com"; E="88 (28)12345678"; F="+84123456887"; G="https://www.google.com.vn/?hl=vi"; H="Gooogle"; I="https://picsum.photos/200/304" },
    @{ Row=7; A=6; B="Nguyen Van B"; C="Head of IT"; D="abcd@gmail.com"; E="89 (28)12345678"; F="+84123456887"; G="https://www.google.com.vn/?hl=vi"; H="Gooogle"; I="https://picsum.photos/200/305" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
    $ws.Range("H$n").Value = $r.H
    $ws.Range("I$n").Value = $r.I
}

# Rows shrink back to the default height now that the long "Managing
# Director..." business title text is gone (row 3 keeps its own custom
# 25.5pt height, unaffected since its text doesn't change length-wise).
$ws.Rows.Item(2).EntireRow.AutoFit()
$ws.Rows.Item(4).EntireRow.AutoFit()
$ws.Rows.Item(5).EntireRow.AutoFit()
$ws.Rows.Item(6).EntireRow.AutoFit()
$ws.Rows.Item(7).EntireRow.AutoFit()

# --- Hyperlinks -------------------------------------------------------
# Remove every existing hyperlink (mail/company/photo links) so they can be
# rebuilt pointing at the new targets; D3:D7 and G3:G7 are now combined
# into single multi-cell hyperlinks instead of five separate ones each.
$existing = @()
foreach ($h in $ws.Hyperlinks) { $existing += $h }
foreach ($h in $existing) { $h.Delete() }

$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:abcd@gmail.com", "", "", "abcd@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), "https://www.google.com.vn/?hl=vi", "", "", "Gooogle") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://picsum.photos/200/300", "", "", "https://picsum.photos/200/300") | Out-Null

$ws.Hyperlinks.Add($ws.Range("I3"), "https://picsum.photos/200/300", "", "", "https://picsum.photos/200/300") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I4"), "https://picsum.photos/200/300", "", "", "https://picsum.photos/200/300") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I5"), "https://picsum.photos/200/300", "", "", "https://picsum.photos/200/300") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I6"), "https://picsum.photos/200/300", "", "", "https://picsum.photos/200/300") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I7"), "https://picsum.photos/200/300", "", "", "https://picsum.photos/200/300") | Out-Null

$ws.Hyperlinks.Add($ws.Range("G3:G7"), "https://www.google.com.vn/?hl=vi", "", "", "https://www.google.com.vn/?hl=vi") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3:D7"), "mailto:abcd@gmail.com", "", "", "abcd@gmail.com") | Out-Null

# --- Column G is much wider now (it shows the full Google URL) ---
$ws.Columns.Item(7).ColumnWidth = 40

# --- Selection cursor moved from E14 to E19 ---
$ws.Range("E19").Select()
